$d = $word.ActiveDocument

$d.Content.Find.Execute("Registrar Propiedades", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Registrar Propiedade", 2)
